# AIP-494 AIP-790 Updated Test Data File
#
# Appends two new data rows (rows 3 and 4) to Sheet1, mirroring the
# structure of the existing template row (row 2):
#   - Row 3: PostFaultTime 300 -> 1000, RecordDuration 500 -> 700,
#            MaxDFR 1000 -> 700 (new "700" value introduced).
#   - Row 4: identical to row 2 except MaxDFR 1000 -> 500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 ----
$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "Admin"
$ws.Range("C3").Value = "IND_DAU_51"
$ws.Range("D3").Value = "IDM+18"
$ws.Range("E3").Value = "10.75.58.51"
$ws.Range("F3").Value = 409026540
$ws.Range("G3").Value = "'100"
$ws.Range("H3").Value = "'1000"
$ws.Range("I3").Value = "'200"
$ws.Range("J3").Value = "'700"
$ws.Range("K3").Value = "CAM_731.seq"
$ws.Range("L3").Value = "'700"
$ws.Range("M3").Value = "RMS H1 - Channel 1"
$ws.Range("N3").Value = "Over"
$ws.Range("O3").Value = "Secondary"
$ws.Range("P3").Value = "'70"

# ---- Row 4 ----
$ws.Range("A4").Value = "Admin"
$ws.Range("B4").Value = "Admin"
$ws.Range("C4").Value = "IND_DAU_51"
$ws.Range("D4").Value = "IDM+18"
$ws.Range("E4").Value = "10.75.58.51"
$ws.Range("F4").Value = 409026540
$ws.Range("G4").Value = "'100"
$ws.Range("H4").Value = "'300"
$ws.Range("I4").Value = "'200"
$ws.Range("J4").Value = "'500"
$ws.Range("K4").Value = "CAM_731.seq"
$ws.Range("L4").Value = "'500"
$ws.Range("M4").Value = "RMS H1 - Channel 1"
$ws.Range("N4").Value = "Over"
$ws.Range("O4").Value = "Secondary"
$ws.Range("P4").Value = "'70"

# Copy row 2's cell formatting (the quote-prefixed / numeric style used
# by columns F,G,H,I,J,L,P) onto the two new rows, without touching the
# values that were just written.
$ws.Range("A2:P2").Copy()
$ws.Range("A3:P3").PasteSpecial(-4122)
$ws.Range("A4:P4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the post-edit selection state (active cell H4).
$ws.Range("H4").Select()
